$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: round all numeric sensor readings down to 2 decimal places
# ("custom accuracy") while the date/time stamp in column A is untouched.
$ws.Range("B5").Value = 14.41
$ws.Range("C5").Value = 10.49
$ws.Range("D5").Value = 1.04
$ws.Range("E5").Value = 31.09
$ws.Range("F5").Value = 25.62
$ws.Range("G5").Value = 11.34
$ws.Range("H5").Value = 41.51
$ws.Range("I5").Value = 17.45
$ws.Range("J5").Value = 7.67
$ws.Range("K5").Value = 11.44
$ws.Range("L5").Value = 12.56
$ws.Range("M5").Value = 13.15
$ws.Range("N5").Value = 3.62
$ws.Range("O5").Value = 11.28
$ws.Range("P5").Value = 15.97
$ws.Range("Q5").Value = 9.59
$ws.Range("R5").Value = 0.8100000000000001
$ws.Range("S5").Value = 0.65
$ws.Range("T5").Value = 164.33
$ws.Range("U5").Value = 31.44
$ws.Range("V5").Value = 10.41
$ws.Range("W5").Value = 21.02
$ws.Range("X5").Value = 11.2
$ws.Range("Y5").Value = 1.69
$ws.Range("Z5").Value = 20.19
$ws.Range("AA5").Value = 9.199999999999999
$ws.Range("AB5").Value = 8.220000000000001
$ws.Range("AC5").Value = 9.65
$ws.Range("AD5").Value = 13.13
$ws.Range("AE5").Value = 0.5600000000000001
$ws.Range("AF5").Value = 37.36
$ws.Range("AG5").Value = 5.81
$ws.Range("AH5").Value = 13.02

# The trailing data point (old row 6) is dropped entirely; the dataset now
# ends at row 5 (used range shrinks from A1:AH6 to A1:AH5).
$ws.Rows.Item(6).Delete()

# Column AH (J33) no longer needs to fit the longer "20.39" string that used
# to live in the deleted row, so its autosized width narrows from 8 to 7
# characters.
$ws.Columns.Item(34).ColumnWidth = 6.14
